$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph.
$verFind = $d.Content
$okVer = $verFind.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$verPara = $verFind.Paragraphs(1)

# The blank paragraph immediately preceding it also needs to be removed.
$emptyBefore = $verPara.Previous(1)

# Locate the "© 2020 . Contact: ..." paragraph that follows.
$copyrightFind = $d.Content
$okCopyright = $copyrightFind.Find.Execute("Contact: luizeleno@usp.br", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyrightPara = $copyrightFind.Paragraphs(1)

# Build a range spanning from the start of the blank paragraph through the end
# of the copyright paragraph (including its paragraph mark) and delete it,
# removing all three paragraphs in one shot.
$delRange = $d.Range($emptyBefore.Range.Start, $copyrightPara.Range.End)
$delRange.Delete()
